$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (age group 14)
$ws.Range("B2").Value = 55.2953815419129
$ws.Range("K2").Value = 56.8474042011213
$ws.Range("L2").Value = 48.4302557813727
$ws.Range("N2").Value = 48.666425997025

# Row 3 (age group 15)
$ws.Range("B3").Value = 40.5493307668479
$ws.Range("K3").Value = 36.9055302021312
$ws.Range("L3").Value = 37.5977317839485
$ws.Range("N3").Value = 41.6894250824717

# Row 4 (age group 16)
$ws.Range("B4").Value = 34.7872415482579
$ws.Range("K4").Value = 30.4357942185907
$ws.Range("N4").Value = 43.0029533260978
